# Add 10 new workout rows (370-379) to Sheet1, mirroring the weekly
# scoreboard update captured in the source diff. Each new row is seeded by
# copying the formatting (incl. the date number format on column B) from
# the row immediately above it, then the actual values are written in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 370
$ws.Range("A369:M369").Copy()
$ws.Range("A370:M370").PasteSpecial(-4122)
$ws.Range("A370").Value2 = "Matt"
$ws.Range("B370").Value2 = 45516
$ws.Range("C370").Value2 = "Workout"
$ws.Range("D370").Value2 = 36
$ws.Range("E370").Value2 = 0
$ws.Range("F370").Value2 = 0
$ws.Range("G370").Value2 = 12
$ws.Range("H370").Value2 = 17
$ws.Range("I370").Value2 = 5
$ws.Range("J370").Value2 = 2
$ws.Range("K370").Value2 = 0
$ws.Range("L370").Value2 = "Agile Antelope"
$ws.Range("M370").Value2 = 10

# Row 371
$ws.Range("A370:M370").Copy()
$ws.Range("A371:M371").PasteSpecial(-4122)
$ws.Range("A371").Value2 = "Steven"
$ws.Range("B371").Value2 = 45516
$ws.Range("C371").Value2 = "Workout"
$ws.Range("D371").Value2 = 25
$ws.Range("E371").Value2 = 0
$ws.Range("F371").Value2 = 0
$ws.Range("G371").Value2 = 21
$ws.Range("H371").Value2 = 3
$ws.Range("I371").Value2 = 1
$ws.Range("J371").Value2 = 0
$ws.Range("K371").Value2 = 0
$ws.Range("L371").Value2 = "Brave Leopard"
$ws.Range("M371").Value2 = 10

# Row 372
$ws.Range("A371:M371").Copy()
$ws.Range("A372:M372").PasteSpecial(-4122)
$ws.Range("A372").Value2 = "Steven"
$ws.Range("B372").Value2 = 45516
$ws.Range("C372").Value2 = "Walk"
$ws.Range("D372").Value2 = 35
$ws.Range("E372").Value2 = 1.72
$ws.Range("F372").Value2 = 105
$ws.Range("G372").Value2 = 35
$ws.Range("H372").Value2 = 0
$ws.Range("I372").Value2 = 0
$ws.Range("J372").Value2 = 0
$ws.Range("K372").Value2 = 0
$ws.Range("L372").Value2 = "Brave Leopard"
$ws.Range("M372").Value2 = 10

# Row 373
$ws.Range("A372:M372").Copy()
$ws.Range("A373:M373").PasteSpecial(-4122)
$ws.Range("A373").Value2 = "Matt"
$ws.Range("B373").Value2 = 45517
$ws.Range("C373").Value2 = "Run"
$ws.Range("D373").Value2 = 27
$ws.Range("E373").Value2 = 3.21
$ws.Range("F373").Value2 = 240
$ws.Range("G373").Value2 = 0
$ws.Range("H373").Value2 = 1
$ws.Range("I373").Value2 = 6
$ws.Range("J373").Value2 = 16
$ws.Range("K373").Value2 = 2
$ws.Range("L373").Value2 = "Agile Antelope"
$ws.Range("M373").Value2 = 10

# Row 374
$ws.Range("A373:M373").Copy()
$ws.Range("A374:M374").PasteSpecial(-4122)
$ws.Range("A374").Value2 = "Matt"
$ws.Range("B374").Value2 = 45517
$ws.Range("C374").Value2 = "Walk"
$ws.Range("D374").Value2 = 1
$ws.Range("E374").Value2 = 0.06
$ws.Range("F374").Value2 = 0
$ws.Range("G374").Value2 = 1
$ws.Range("H374").Value2 = 0
$ws.Range("I374").Value2 = 0
$ws.Range("J374").Value2 = 0
$ws.Range("K374").Value2 = 0
$ws.Range("L374").Value2 = "Agile Antelope"
$ws.Range("M374").Value2 = 10

# Row 375
$ws.Range("A374:M374").Copy()
$ws.Range("A375:M375").PasteSpecial(-4122)
$ws.Range("A375").Value2 = "Steven"
$ws.Range("B375").Value2 = 45517
$ws.Range("C375").Value2 = "Workout"
$ws.Range("D375").Value2 = 15
$ws.Range("E375").Value2 = 0
$ws.Range("F375").Value2 = 0
$ws.Range("G375").Value2 = 15
$ws.Range("H375").Value2 = 0
$ws.Range("I375").Value2 = 0
$ws.Range("J375").Value2 = 0
$ws.Range("K375").Value2 = 0
$ws.Range("L375").Value2 = "Brave Leopard"
$ws.Range("M375").Value2 = 10

# Row 376
$ws.Range("A375:M375").Copy()
$ws.Range("A376:M376").PasteSpecial(-4122)
$ws.Range("A376").Value2 = "Steven"
$ws.Range("B376").Value2 = 45517
$ws.Range("C376").Value2 = "Walk"
$ws.Range("D376").Value2 = 18
$ws.Range("E376").Value2 = 0.96
$ws.Range("F376").Value2 = 98
$ws.Range("G376").Value2 = 15
$ws.Range("H376").Value2 = 3
$ws.Range("I376").Value2 = 0
$ws.Range("J376").Value2 = 0
$ws.Range("K376").Value2 = 0
$ws.Range("L376").Value2 = "Brave Leopard"
$ws.Range("M376").Value2 = 10

# Row 377
$ws.Range("A376:M376").Copy()
$ws.Range("A377:M377").PasteSpecial(-4122)
$ws.Range("A377").Value2 = "Eric"
$ws.Range("B377").Value2 = 45517
$ws.Range("C377").Value2 = "Run"
$ws.Range("D377").Value2 = 44
$ws.Range("E377").Value2 = 3.86
$ws.Range("F377").Value2 = 131
$ws.Range("G377").Value2 = 0
$ws.Range("H377").Value2 = 14
$ws.Range("I377").Value2 = 24
$ws.Range("J377").Value2 = 0
$ws.Range("K377").Value2 = 0
$ws.Range("L377").Value2 = "Sauntering Hippo"
$ws.Range("M377").Value2 = 10

# Row 378
$ws.Range("A377:M377").Copy()
$ws.Range("A378:M378").PasteSpecial(-4122)
$ws.Range("A378").Value2 = "Steven"
$ws.Range("B378").Value2 = 45517
$ws.Range("C378").Value2 = "Workout"
$ws.Range("D378").Value2 = 33
$ws.Range("E378").Value2 = 0
$ws.Range("F378").Value2 = 0
$ws.Range("G378").Value2 = 14
$ws.Range("H378").Value2 = 17
$ws.Range("I378").Value2 = 2
$ws.Range("J378").Value2 = 0
$ws.Range("K378").Value2 = 0
$ws.Range("L378").Value2 = "Brave Leopard"
$ws.Range("M378").Value2 = 10

# Row 379
$ws.Range("A378:M378").Copy()
$ws.Range("A379:M379").PasteSpecial(-4122)
$ws.Range("A379").Value2 = "Steven"
$ws.Range("B379").Value2 = 45517
$ws.Range("C379").Value2 = "Walk"
$ws.Range("D379").Value2 = 35
$ws.Range("E379").Value2 = 1.76
$ws.Range("F379").Value2 = 157
$ws.Range("G379").Value2 = 34
$ws.Range("H379").Value2 = 1
$ws.Range("I379").Value2 = 0
$ws.Range("J379").Value2 = 0
$ws.Range("K379").Value2 = 0
$ws.Range("L379").Value2 = "Brave Leopard"
$ws.Range("M379").Value2 = 10

# Move the selection to the new last cell, matching where Excel leaves the
# cursor after appending rows at the bottom of the table.
$ws.Range("M379").Select() | Out-Null
